$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column E (Approved/Send Confirmation Email
# and the trailing blank columns all shift right by two: E->G, F->H, ... J->L)
$ws.Range("E1:F1").EntireColumn.Insert()

# Match the width of the other name/email columns (B:D) for the two new columns
$ws.Range("E1:F1").EntireColumn.ColumnWidth = $ws.Range("D1").EntireColumn.ColumnWidth

# New column headers
$ws.Range("E1").Value = "Phone"
$ws.Range("F1").Value = "WhatsApp Enabled"

# New "Phone" values for the two data rows
$ws.Range("E2").Value = 9999999999
$ws.Range("E3").Value = 8888888888

# New "WhatsApp Enabled" values mirror the (now shifted) "Approved" column
$ws.Range("F2").Value = $ws.Range("G2").Value2
$ws.Range("F3").Value = $ws.Range("G3").Value2

# Restore the active cell selection to F4
$ws.Range("F4").Select()
